$wb = $excel.ActiveWorkbook

$longMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1700165b9a9c9b35861086596499a724c2d705c/e2e/55d63824-075b-4635-a5ee-0e8921f76cd9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3f4e444481152b88c6954feb444f5b0dbcffcc9/e2e/55d63824-075b-4635-a5ee-0e8921f76cd9.md."

# zh-cn sheet: row 7 is the handback report row for 55d63824-075b-4635-a5ee-0e8921f76cd9
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("J7").Value = "55d63824-075b-4635-a5ee-0e8921f76cd9.2b327231f3a4f2ba38baae7e2e5829f934d18aac.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-19 06:53:47"
$wsZh.Range("P7").Value = $longMessage
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b3f4e444481152b88c6954feb444f5b0dbcffcc9/e2e/55d63824-075b-4635-a5ee-0e8921f76cd9.md", $null, $null, "55d63824-075b-4635-a5ee-0e8921f76cd9.md")

# de-de sheet: row 7 is the handback report row for 55d63824-075b-4635-a5ee-0e8921f76cd9
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("J7").Value = "55d63824-075b-4635-a5ee-0e8921f76cd9.2b327231f3a4f2ba38baae7e2e5829f934d18aac.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-19 06:53:54"
$wsDe.Range("P7").Value = $longMessage
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b3f4e444481152b88c6954feb444f5b0dbcffcc9/e2e/55d63824-075b-4635-a5ee-0e8921f76cd9.md", $null, $null, "55d63824-075b-4635-a5ee-0e8921f76cd9.md")
